$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pinouts")
$ws.Range("W3").Value = "v1"
$rng = $ws.Range("W3")
try {
  $rng.Format.Fill.Color = "#4BACC6"
  Write-Host "set format.fill.color ok"
} catch {
  Write-Host "ERR1:" $_
}
try {
  $rng.Format.Fill.TintAndShade = -0.9
  Write-Host "set format.fill.tint ok"
} catch {
  Write-Host "ERR2:" $_
}
Write-Host "Color:" $rng.Interior.Color
